$wb = $excel.ActiveWorkbook

# --- Sheet1: lot_management_rough_exmple - remove the trailing blank rows 16:19 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A16:G19").Delete()

# --- Sheet2: VBL lot profiling is done ---
$ws2 = $wb.Worksheets.Item("Sheet2")

# E1 header cell: normalize its format (was using a stray duplicate font/style)
$e1 = $ws2.Cells.Item(1, 5)
$e1.HorizontalAlignment = 1
$e1.VerticalAlignment = -4107
$e1.WrapText = $false
$e1.NumberFormat = "General"

# E11: replace the hard-coded average price with the running-average formula
# used throughout the rest of column E (this ripples into G11/I11/K11).
$ws2.Range("E11").Formula = "=(G10+C11+D11)/F11"

# Move the selection/active cell to reflect where work left off.
$ws2.Activate()
$ws2.Range("E10").Select()
